$wb = $excel.ActiveWorkbook

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2099.7778
$ws.Range("I51").Value = 1500
$ws.Range("J51").Value = 2174.75
$ws.Range("K51").Value = 1500
$ws.Range("L51").Value = 2174.75
$ws.Range("M51").Value = -1016
$ws.Range("N51").Value = -3142.75

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4700
$ws.Range("I64").Value = 5333.3335
$ws.Range("K64").Value = 5333.3335
$ws.Range("M64").Value = -5085.3335

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4700
$ws.Range("I67").Value = 5333.3335
$ws.Range("K67").Value = 5333.3335
$ws.Range("M67").Value = -4475.3335

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4350780.5
$ws.Range("I76").Value = 4548407
$ws.Range("K76").Value = 4548407
$ws.Range("M76").Value = -4548092

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4350780.5
$ws.Range("I79").Value = 4548407
$ws.Range("K79").Value = 4548407
$ws.Range("M79").Value = -4547315

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1594.25
$ws.Range("I111").Value = 1859
$ws.Range("J111").Value = 800
$ws.Range("K111").Value = 5577
$ws.Range("L111").Value = 2400
$ws.Range("M111").Value = -2510
$ws.Range("N111").Value = -8534

# ALC row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 34942.5
$ws.Range("J123").Value = 34942.5
$ws.Range("L123").Value = 34942.5
$ws.Range("N123").Value = -44742.5

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1128.4546
$ws.Range("J137").Value = 1408.25
$ws.Range("L137").Value = 4224.75
$ws.Range("N137").Value = -9324.75

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10755.013
$ws.Range("I32").Value = 9153.197
$ws.Range("J32").Value = 27002
$ws.Range("K32").Value = 9153.197
$ws.Range("L32").Value = 27002
$ws.Range("M32").Value = -8866.197
$ws.Range("N32").Value = -27576

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1881.1143
$ws.Range("I63").Value = 1877.6177
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 1877.6177
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1191.6177
$ws.Range("N63").Value = -3372

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1881.1143
$ws.Range("I66").Value = 1877.6177
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 9388.0885
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -5956.0885
$ws.Range("N66").Value = -16864

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1086.2
$ws.Range("I122").Value = 952.41174
$ws.Range("J122").Value = 1370.5
$ws.Range("K122").Value = 2857.23522
$ws.Range("L122").Value = 4111.5
$ws.Range("M122").Value = -407.23522
$ws.Range("N122").Value = -9011.5

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4623.0293
$ws.Range("I105").Value = 4014.3928
$ws.Range("J105").Value = 7463.3335
$ws.Range("K105").Value = 4014.3928
$ws.Range("L105").Value = 7463.3335
$ws.Range("M105").Value = -2267.3928
$ws.Range("N105").Value = -10957.3335

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3088.804
$ws.Range("I31").Value = 2952.6487
$ws.Range("J31").Value = 3448.6428
$ws.Range("K31").Value = 2952.6487
$ws.Range("L31").Value = 3448.6428
$ws.Range("M31").Value = -2657.6487
$ws.Range("N31").Value = -4038.6428

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3088.804
$ws.Range("I34").Value = 2952.6487
$ws.Range("J34").Value = 3448.6428
$ws.Range("K34").Value = 2952.6487
$ws.Range("L34").Value = 3448.6428
$ws.Range("M34").Value = -2750.6487
$ws.Range("N34").Value = -3852.6428

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4528.4614
$ws.Range("I62").Value = 5607.143
$ws.Range("J62").Value = 3270
$ws.Range("K62").Value = 5607.143
$ws.Range("L62").Value = 3270
$ws.Range("M62").Value = -4983.143
$ws.Range("N62").Value = -4518

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4528.4614
$ws.Range("I65").Value = 5607.143
$ws.Range("J65").Value = 3270
$ws.Range("K65").Value = 28035.715
$ws.Range("L65").Value = 16350
$ws.Range("M65").Value = -24915.715
$ws.Range("N65").Value = -22590

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 5721.375
$ws.Range("I94").Value = 4522.4
$ws.Range("J94").Value = 6266.364
$ws.Range("K94").Value = 4522.4
$ws.Range("L94").Value = 6266.364
$ws.Range("M94").Value = -4071.4
$ws.Range("N94").Value = -7168.364

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1859.826
$ws.Range("I132").Value = 1185.3572
$ws.Range("K132").Value = 3556.0716
$ws.Range("M132").Value = -1026.0716

# CUL row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3433.7693
$ws.Range("I58").Value = 1699.5
$ws.Range("J58").Value = 3749.0908
$ws.Range("K58").Value = 5098.5
$ws.Range("L58").Value = 11247.2724
$ws.Range("M58").Value = -4970.5
$ws.Range("N58").Value = -11503.2724

# CUL row 118
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1427.25
$ws.Range("I118").Value = 354.5
$ws.Range("J118").Value = 2500
$ws.Range("K118").Value = 1063.5
$ws.Range("L118").Value = 7500
$ws.Range("M118").Value = 179.5
$ws.Range("N118").Value = -9986

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4508.3335
$ws.Range("I70").Value = 4566.619
$ws.Range("J70").Value = 4304.3335
$ws.Range("K70").Value = 4566.619
$ws.Range("L70").Value = 4304.3335
$ws.Range("M70").Value = -4296.619
$ws.Range("N70").Value = -4844.3335

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4508.3335
$ws.Range("I73").Value = 4566.619
$ws.Range("J73").Value = 4304.3335
$ws.Range("K73").Value = 4566.619
$ws.Range("L73").Value = 4304.3335
$ws.Range("M73").Value = -3630.619
$ws.Range("N73").Value = -6176.3335

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3994.75
$ws.Range("I80").Value = 3958.1667
$ws.Range("J80").Value = 4049.625
$ws.Range("K80").Value = 3958.1667
$ws.Range("L80").Value = 4049.625
$ws.Range("M80").Value = -2960.1667
$ws.Range("N80").Value = -6045.625

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3994.75
$ws.Range("I83").Value = 3958.1667
$ws.Range("J83").Value = 4049.625
$ws.Range("K83").Value = 19790.8335
$ws.Range("L83").Value = 20248.125
$ws.Range("M83").Value = -14798.8335
$ws.Range("N83").Value = -30232.125

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1645706.8
$ws.Range("I122").Value = 2632470.8
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 7897412.399999999
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -7894962.399999999
$ws.Range("N122").Value = -8200

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1858.8948
$ws.Range("I107").Value = 1825.8235
$ws.Range("J107").Value = 2140
$ws.Range("K107").Value = 5477.470499999999
$ws.Range("L107").Value = 6420
$ws.Range("M107").Value = -3557.470499999999
$ws.Range("N107").Value = -10260
